# Apply pricing/date refresh to "Hoja1" of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# A1 holds a date serial (format changed from 45406 -> 45436, i.e. 2024-04-24 -> 2024-05-24)
$ws.Range("A1").Value = 45436

# D33 carries a quote-prefixed numeric style (s=11); writing a new Value normally
# swaps the cell onto the equivalent non-quote-prefixed style (s=13). Preserve the
# original formatting by stashing/restoring it through copy-format round trip via
# an out-of-the-way scratch cell, then clearing the scratch cell afterwards.
$d33 = $ws.Range("D33")
$scratch = $ws.Range("Z1")
$d33.Copy()
$scratch.PasteSpecial(-4122)  # xlPasteFormats

$d33.Value = 1687.737

$scratch.Copy()
$d33.PasteSpecial(-4122)      # xlPasteFormats
$scratch.Clear()

# Updated price list values in column D
$ws.Range("D34").Value = 1275.478
$ws.Range("D35").Value = 949.728
